$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-06-15"

$ws.Range("A7").Value = "June (through 06-15)"

$ws.Range("D7").Value = 28
$ws.Range("E7").Value = 31
$ws.Range("F7").Value = 23
$ws.Range("G7").Value = 57
$ws.Range("H7").Value = 53
$ws.Range("I7").Value = 69

$ws.Range("D8").Value = 344
$ws.Range("E8").Value = 326
$ws.Range("F8").Value = 227
$ws.Range("G8").Value = 415
$ws.Range("H8").Value = 684
$ws.Range("I8").Value = 732
